$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.314.36'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').Value = '3.392.93'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.73'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.78'
$ws.Range('E6').Value = '  +1.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('E9').Value = '  +8.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.591'
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.77'
$ws.Range('E11').Value = '  +4.13%  '
$ws.Range('E12').Value = '  +3.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '681.15'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('D15').Value = '3.938.93'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = '69.383.47'
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('D17').Value = '3.395.88'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.75'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.39'
$ws.Range('E20').Value = '  +3.04%  '
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.43'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.16'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.61'
$ws.Range('E24').Value = '  +4.98%  '
$ws.Range('E25').Value = '  +1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.74'
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.65'
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.30'
$ws.Range('E28').Value = '  +3.89%  '
$ws.Range('E29').Value = '  +1.98%  '
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('E31').Value = '  +11.99%  '
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '555.56'
$ws.Range('E33').Value = '  -3.29%  '
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '3.709.10'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.75'
$ws.Range('E38').Value = '  +3.28%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.140'
$ws.Range('E39').Value = '  +6.37%  '
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('E41').Value = '  +4.64%  '
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0424'
$ws.Range('E44').Value = '  +4.24%  '
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('E48').Value = '  +5.65%  '
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '131.68'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.58'
$ws.Range('E51').Value = '  -2.25%  '
